$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = 5
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = 8
$ws.Range("I28").Value = 9

$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 2
$ws.Range("E32").Value = 3
$ws.Range("F32").Value = 4

$ws.Range("C36").Value = 4
$ws.Range("D36").Value = 2
$ws.Range("E36").Value = 1
$ws.Range("G36").Value = -3

$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("F36").Select()
